# Actualización automática 2025-07-16 10:50:09
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M3").Value = 3317.32
$wsGrupo.Range("M24").Value = "2 de 22"

# --- Sheet "VENTA MENSUAL" ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F3").Value = 3317.32
$wsMensual.Range("F24").Value = 9584.200000000001

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumpl.Range("D16").Value = 9584.200000000001
$wsCumpl.Range("E16").Value = 29172.34
$wsCumpl.Range("F16").Value = 0.2472924569633925

$wsCumpl.Range("D19").Value = 9584.200000000001
$wsCumpl.Range("E19").Value = 48638.80386304604
$wsCumpl.Range("F19").Value = 0.1646119122013054

$wsCumpl.Columns.Item(4).ColumnWidth = 11.166666666666666
